$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item('展览')
$ws.Range('F4').Value = 663
$ws.Range('F5').Value = 2967
$ws.Range('F7').Value = 243
$ws.Range('F10').Value = 6992
$ws.Range('F12').Value = 108
$ws.Range('F13').Value = 367
$ws.Range('F14').Value = 612
$ws.Range('F15').Value = 1510
$ws.Range('F16').Value = 1128
$ws.Range('F17').Value = 2264
$ws.Range('F18').Value = 1514
$ws.Range('F20').Value = 1122
$ws.Range('F22').Value = 2
$ws.Range('F23').Value = 194
$ws.Range('F25').Value = 19
$ws.Range('F26').Value = 1762
$ws.Range('F27').Value = 1706
$ws.Range('F28').Value = 1037
$ws.Range('F29').Value = 39
$ws.Range('F30').Value = 1678
$ws.Range('F31').Value = 1241
$ws.Range('F32').Value = 144
$ws.Range('F34').Value = 6
$ws.Range('F35').Value = 1066
$ws.Range('F36').Value = 441
$ws.Range('F37').Value = 27
$ws.Range('F38').Value = 2515
$ws.Range('F39').Value = 2757
$ws.Range('F41').Value = 28
$ws.Range('F42').Value = 189
$ws.Range('B43').Value = '2024-09-15'
$ws.Range('C43').Value = '上海·第十一届次元鹿角动漫游戏展·月映中秋'
$ws.Range('D43').Value = '曹杨路1888号 复悦荟'
$ws.Range('E43').Value = '2024.09.15 10:00-09.17 17:00'
$ws.Range('F43').Value = 1
$ws.Range('G43').Value = 65
$ws.Range('H43').Value = 'https://show.bilibili.com/platform/detail.html?id=90799'
$ws.Range('I43').Value = '//i1.hdslb.com/bfs/openplatform/202408/ZJ1OBaIK1723635427833.png'
$ws.Range('C44').Value = '上海·LookLook动漫嘉年华3th'
$ws.Range('D44').Value = '曹安公路4218号 上海国际短视频中心'
$ws.Range('E44').Value = '2024.09.21 10:00-09.22 17:30'
$ws.Range('F44').Value = 20
$ws.Range('G44').Value = 68
$ws.Range('H44').Value = 'https://show.bilibili.com/platform/detail.html?id=90495'
$ws.Range('I44').Value = '//i0.hdslb.com/bfs/openplatform/202408/8gEx4nQa1723033308040.jpeg'
$ws.Range('C45').Value = '上海·代号鸢同人only之谁偷了我的狗'
$ws.Range('D45').Value = '申旺路18号3幢 丽丽薇里婚嫁礼堂婚宴(闵行店)'
$ws.Range('E45').Value = '2024.09.21 10:00-09.21 20:00'
$ws.Range('F45').Value = 31
$ws.Range('G45').Value = 69
$ws.Range('H45').Value = 'https://show.bilibili.com/platform/detail.html?id=90490'
$ws.Range('I45').Value = '//i2.hdslb.com/bfs/openplatform/202408/oBLHzc1O1723037740163.jpeg'
$ws.Range('B46').Value = '2024-09-21'
$ws.Range('C46').Value = '上海·咒术回战同人ONLY'
$ws.Range('E46').Value = '2024.09.21 10:00-09.21 17:00'
$ws.Range('F46').Value = 329
$ws.Range('H46').Value = 'https://show.bilibili.com/platform/detail.html?id=89162'
$ws.Range('I46').Value = '//i0.hdslb.com/bfs/openplatform/202407/jx24gc5U1720746534824.jpeg'
$ws.Range('B47').Value = '2024-09-22'
$ws.Range('C47').Value = '上海·火影同人only'
$ws.Range('E47').Value = '2024.09.22 10:00-09.22 17:00'
$ws.Range('F47').Value = 126
$ws.Range('G47').Value = 65
$ws.Range('H47').Value = 'https://show.bilibili.com/platform/detail.html?id=89551'
$ws.Range('I47').Value = '//i2.hdslb.com/bfs/openplatform/202407/DBYRHEuC1721123198969.jpeg'
$ws.Range('C48').Value = '上海·原神X星穹铁道x绝区零同人ONLY'
$ws.Range('D48').Value = '顾村镇蕰川路6号 智慧湾科创园'
$ws.Range('E48').Value = '2024.10.01 10:00-10.02 17:00'
$ws.Range('F48').Value = 173
$ws.Range('G48').Value = 60
$ws.Range('H48').Value = 'https://show.bilibili.com/platform/detail.html?id=90135'
$ws.Range('I48').Value = '//i0.hdslb.com/bfs/openplatform/202407/FF8HGnt01722418798545.jpeg'

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item('演出')
$ws.Range('F7').Value = 175
$ws.Range('F10').Value = 32
$ws.Range('F12').Value = 186
$ws.Range('F17').Value = 164
$ws.Range('F19').Value = 47
$ws.Range('F20').Value = 61
$ws.Range('F23').Value = 485

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item('本地生活')
$ws.Range('F4').Value = 545
$ws.Range('F6').Value = 1706
$ws.Range('F7').Value = 1857
$ws.Range('F8').Value = 2759
$ws.Range('F9').Value = 1035
$ws.Range('F10').Value = 953
$ws.Range('F12').Value = 296
$ws.Range('F13').Value = 1528
$ws.Range('F14').Value = 7425

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item('全部类型')
$ws.Range('F3').Value = 663
$ws.Range('F4').Value = 2967
$ws.Range('F5').Value = 243
$ws.Range('F6').Value = 1706
$ws.Range('F7').Value = 2759
$ws.Range('F8').Value = 6992
$ws.Range('F9').Value = 1035
$ws.Range('F11').Value = 108
$ws.Range('F12').Value = 367
$ws.Range('F13').Value = 175
$ws.Range('F14').Value = 1528
$ws.Range('F15').Value = 612
$ws.Range('F16').Value = 1510
$ws.Range('F17').Value = 1128
$ws.Range('F18').Value = 2264
$ws.Range('F19').Value = 1514
$ws.Range('F21').Value = 186
$ws.Range('F22').Value = 1122
$ws.Range('F25').Value = 19
$ws.Range('F26').Value = 1762
$ws.Range('F27').Value = 1037
$ws.Range('F28').Value = 39
$ws.Range('F29').Value = 1678
$ws.Range('F30').Value = 1241
$ws.Range('F31').Value = 144
$ws.Range('F33').Value = 1067
$ws.Range('F34').Value = 61
$ws.Range('F36').Value = 485
$ws.Range('F37').Value = 441
$ws.Range('F38').Value = 27
$ws.Range('F39').Value = 2515
$ws.Range('F40').Value = 2757
$ws.Range('F42').Value = 189
$ws.Range('F43').Value = 31
$ws.Range('F44').Value = 329
$ws.Range('F46').Value = 173
